# Update countries & provincias Spain
# Refresh of the COVID-19 country data table: some countries' case counts
# were updated, which changes their rank (the sheet is sorted descending by
# "Casos totales" / column B), so a few adjacent rows swap which country
# they show.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- China (row 10) : updated figures, no reordering ---
$ws.Range("B10").Value = 82367
$ws.Range("C10").Value = 26
$ws.Range("D10").Value = 77944
$ws.Range("E10").Value = 1081
$ws.Range("F10").Value = 89

# --- Corea del Sur (row 26) : updated figures, no reordering ---
$ws.Range("B26").Value = 10635
$ws.Range("C26").Value = 22
$ws.Range("D26").Value = 7829
$ws.Range("E26").Value = 2576
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 230

# --- Rows 88/89 : Principado de Andorra overtakes Letonia ---
$ws.Range("A88").Value = "Principado de Andorra"
$ws.Range("B88").Value = 682
$ws.Range("C88").Value = 9
$ws.Range("D88").Value = 169
$ws.Range("E88").Value = 480
$ws.Range("F88").Value = 17
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 33

$ws.Range("A89").Value = "Letonia"
$ws.Range("B89").Value = 675
$ws.Range("C89").Value = 0
$ws.Range("D89").Value = 57
$ws.Range("E89").Value = 613
$ws.Range("F89").Value = 3
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 5

# --- Rows 100/101 : Bolivia overtakes Nigeria ---
$ws.Range("A100").Value = "Bolivia"
$ws.Range("B100").Value = 465
$ws.Range("C100").Value = 24
$ws.Range("D100").Value = 26
$ws.Range("E100").Value = 408
$ws.Range("F100").Value = 3
$ws.Range("G100").Value = 2
$ws.Range("H100").Value = 31

$ws.Range("A101").Value = "Nigeria"
$ws.Range("B101").Value = 442
$ws.Range("C101").Value = 0
$ws.Range("D101").Value = 152
$ws.Range("E101").Value = 277
$ws.Range("F101").Value = 2
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 13

# --- Rows 120-123 : Guatemala & Paraguay overtake Venezuela & Islas Feroe ---
$ws.Range("A120").Value = "Guatemala"
$ws.Range("B120").Value = 214
$ws.Range("C120").Value = 18
$ws.Range("D120").Value = 21
$ws.Range("E120").Value = 186
$ws.Range("F120").Value = 3
$ws.Range("G120").Value = 2
$ws.Range("H120").Value = 7

$ws.Range("A121").Value = "Venezuela"
$ws.Range("B121").Value = 204
$ws.Range("C121").Value = 0
$ws.Range("D121").Value = 111
$ws.Range("E121").Value = 84
$ws.Range("F121").Value = 4
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 9

$ws.Range("A122").Value = "Paraguay"
$ws.Range("B122").Value = 199
$ws.Range("C122").Value = 25
$ws.Range("D122").Value = 30
$ws.Range("E122").Value = 161
$ws.Range("F122").Value = 1
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 8

$ws.Range("A123").Value = "Islas Feroe"
$ws.Range("B123").Value = 184
$ws.Range("C123").Value = 0
$ws.Range("D123").Value = 169
$ws.Range("E123").Value = 15
$ws.Range("F123").Value = 0
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 0
